# Auto-generated Excel COM-interop edit script
# Updates cryptos list (price + volume% columns) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.388.65'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '2.016.82'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '260.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '56.51'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -5.89%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.388'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('E10').Value = '  -4.00%  '
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.35'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -5.21%  '
$ws.Range('D13').Value = '2.314.61'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.807'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.06'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -6.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.27'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.11%  '
$ws.Range('D17').Value = '2.023.44'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '37.368.53'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.98'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').Value = '0.0₃0841'
$ws.Range('E20').Value = '  -2.78%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.20'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.21'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.66'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +7.41%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.36'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.91'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.132'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.00%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.33'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.69%  '
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.40'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.28'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.27'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.05'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.17%  '
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('E42').Value = '  -4.96%  '
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').Value = '1.405.23'
$ws.Range('E44').Value = '  +2.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.35'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '15.83'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.81%  '
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('E48').Value = '  -2.88%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.90'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.18%  '
$ws.Range('D50').Value = '2.206.19'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('E51').Value = '  -6.50%  '
